{"js": "// Docx writer: Use different style for block quotes in notes.\n// Add a new \"Footnote Block Text\" paragraph style, based on (and\n// followed by) \"Footnote Text\", so footnote block quotes can be\n// styled differently from regular footnote text.\n\n// Create the style (minimal stub); properties must be set on a\n// fresh lookup after a sync, otherwise they don't persist.\ncontext.document.addStyle(\"Footnote Block Text\", \"Paragraph\");\nawait context.sync();\n\nconst style = context.document.getStyles().getByNameOrNullObject(\"Footnote Block Text\");\nstyle.baseStyle = \"Footnote Text\";\nstyle.nextParagraphStyle = \"Footnote Text\";\nstyle.priority = 9;\nstyle.unhideWhenUsed = true;\nstyle.quickStyle = true;\n\n// Paragraph formatting: spacing before/after = 100 twips (5 pt),\n// left/right indent = 480 twips (24 pt), no first-line indent.\nconst pf = style.paragraphFormat;\npf.spaceBefore = 5;\npf.spaceAfter = 5;\npf.firstLineIndent = 0;\npf.leftIndent = 24;\npf.rightIndent = 24;\n\nawait context.sync();\n", "ps1": "# Docx writer: Use different style for block quotes in notes.\n# Add a new \"Footnote Block Text\" paragraph style, based on (and\n# followed by) \"Footnote Text\", so footnote block quotes can be\n# styled differently from regular footnote text.\n\n$d = $word.ActiveDocument\n\n# wdStyleTypeParagraph = 1\n$style = $d.Styles.Add(\"Footnote Block Text\", 1)\n\n$style.BaseStyle = \"Footnote Text\"\n$style.NextParagraphStyle = \"Footnote Text\"\n$style.Priority = 9\n$style.UnhideWhenUsed = $true\n$style.QuickStyle = $true\n\n# Paragraph formatting: spacing before/after = 100 twips (5 pt),\n# left/right indent = 480 twips (24 pt), no first-line indent.\n$style.ParagraphFormat.SpaceBefore = 5\n$style.ParagraphFormat.SpaceAfter = 5\n$style.ParagraphFormat.FirstLineIndent = 0\n$style.ParagraphFormat.LeftIndent = 24\n$style.ParagraphFormat.RightIndent = 24\n"}
